$wb = $excel.ActiveWorkbook

# --- Sheet1: per-instance summary table ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Cells.Item(2, 2).Value = -1737.6437622303326
$ws1.Cells.Item(2, 3).Value = 494.973067135
$ws1.Cells.Item(2, 6).Value = 100
$ws1.Cells.Item(2, 7).Value = 110400
$ws1.Cells.Item(2, 8).Value = 120500
$ws1.Cells.Item(2, 9).Value = 10000

$ws1.Cells.Item(3, 2).Value = -1740.0765203442552
$ws1.Cells.Item(3, 3).Value = 1815.573542726
$ws1.Cells.Item(3, 6).Value = 100
$ws1.Cells.Item(3, 7).Value = 110400
$ws1.Cells.Item(3, 8).Value = 120500
$ws1.Cells.Item(3, 9).Value = 10000

$ws1.Cells.Item(4, 2).Value = -1740.3063680164755
$ws1.Cells.Item(4, 3).Value = 631.482673856
$ws1.Cells.Item(4, 6).Value = 100
$ws1.Cells.Item(4, 7).Value = 110400
$ws1.Cells.Item(4, 8).Value = 120500
$ws1.Cells.Item(4, 9).Value = 10000

$ws1.Cells.Item(5, 2).Value = -1735.9954971019556
$ws1.Cells.Item(5, 3).Value = 699.179153258
$ws1.Cells.Item(5, 6).Value = 100
$ws1.Cells.Item(5, 7).Value = 110400
$ws1.Cells.Item(5, 8).Value = 120500
$ws1.Cells.Item(5, 9).Value = 10000

$ws1.Cells.Item(6, 2).Value = -1744.89620665423
$ws1.Cells.Item(6, 3).Value = 519.244673681
$ws1.Cells.Item(6, 6).Value = 100
$ws1.Cells.Item(6, 7).Value = 110400
$ws1.Cells.Item(6, 8).Value = 120500
$ws1.Cells.Item(6, 9).Value = 10000

$ws1.Cells.Item(7, 2).Value = -1725.3743202503338
$ws1.Cells.Item(7, 3).Value = 631.764736205
$ws1.Cells.Item(7, 6).Value = 100
$ws1.Cells.Item(7, 7).Value = 110400
$ws1.Cells.Item(7, 8).Value = 120500
$ws1.Cells.Item(7, 9).Value = 10000

$ws1.Cells.Item(8, 2).Value = -1739.2480133361746
$ws1.Cells.Item(8, 3).Value = 558.887804942
$ws1.Cells.Item(8, 6).Value = 100
$ws1.Cells.Item(8, 7).Value = 110400
$ws1.Cells.Item(8, 8).Value = 120500
$ws1.Cells.Item(8, 9).Value = 10000

$ws1.Cells.Item(9, 2).Value = -1723.387446692379
$ws1.Cells.Item(9, 3).Value = 716.210385767
$ws1.Cells.Item(9, 6).Value = 100
$ws1.Cells.Item(9, 7).Value = 110400
$ws1.Cells.Item(9, 8).Value = 120500
$ws1.Cells.Item(9, 9).Value = 10000

$ws1.Cells.Item(10, 2).Value = -1733.6804531563703
$ws1.Cells.Item(10, 3).Value = 624.165352763
$ws1.Cells.Item(10, 6).Value = 100
$ws1.Cells.Item(10, 7).Value = 110400
$ws1.Cells.Item(10, 8).Value = 120500
$ws1.Cells.Item(10, 9).Value = 10000

$ws1.Cells.Item(11, 2).Value = -1746.9903591802204
$ws1.Cells.Item(11, 3).Value = 969.963798802
$ws1.Cells.Item(11, 6).Value = 100
$ws1.Cells.Item(11, 7).Value = 110400
$ws1.Cells.Item(11, 8).Value = 120500
$ws1.Cells.Item(11, 9).Value = 10000

# --- Per-iteration MP tables (sheets "1".."10") ---
$ws = $wb.Worksheets.Item("1")
$ws.Cells.Item(2, 4).Value = 0.9486819596269531
$ws.Cells.Item(2, 5).Value = 190.71448
$ws.Cells.Item(3, 2).Value = -1737.6437622303326
$ws.Cells.Item(3, 3).Value = 0.07531018375979386
$ws.Cells.Item(3, 4).Value = 479.960756343875

$ws = $wb.Worksheets.Item("2")
$ws.Cells.Item(2, 4).Value = 0.07435619636791992
$ws.Cells.Item(2, 5).Value = 192.02891
$ws.Cells.Item(3, 2).Value = -1740.0765203442552
$ws.Cells.Item(3, 3).Value = 0.0796673759901016
$ws.Cells.Item(3, 4).Value = 1809.3381828707695

$ws = $wb.Worksheets.Item("3")
$ws.Cells.Item(2, 4).Value = 0.0870287155102539
$ws.Cells.Item(2, 5).Value = 191.23263
$ws.Cells.Item(3, 2).Value = -1740.3063680164755
$ws.Cells.Item(3, 3).Value = 0.04341400505932843
$ws.Cells.Item(3, 4).Value = 625.0491888040209

$ws = $wb.Worksheets.Item("4")
$ws.Cells.Item(2, 4).Value = 0.11628421490698242
$ws.Cells.Item(2, 5).Value = 192.98878
$ws.Cells.Item(3, 2).Value = -1735.9954971019556
$ws.Cells.Item(3, 3).Value = 0.08037126701212952
$ws.Cells.Item(3, 4).Value = 692.3040742149391

$ws = $wb.Worksheets.Item("5")
$ws.Cells.Item(2, 4).Value = 0.09610155913867187
$ws.Cells.Item(2, 5).Value = 194.45357
$ws.Cells.Item(3, 2).Value = -1744.89620665423
$ws.Cells.Item(3, 3).Value = 0.09660409261779686
$ws.Cells.Item(3, 4).Value = 513.072218979344

$ws = $wb.Worksheets.Item("6")
$ws.Cells.Item(2, 4).Value = 0.08517435643371582
$ws.Cells.Item(2, 5).Value = 191.51613
$ws.Cells.Item(3, 2).Value = -1725.3743202503338
$ws.Cells.Item(3, 3).Value = 0.07352645462302244
$ws.Cells.Item(3, 4).Value = 625.4828875568791

$ws = $wb.Worksheets.Item("7")
$ws.Cells.Item(2, 4).Value = 0.1123842785354004
$ws.Cells.Item(2, 5).Value = 191.48952
$ws.Cells.Item(3, 2).Value = -1739.2480133361746
$ws.Cells.Item(3, 3).Value = 0.013677691289950836
$ws.Cells.Item(3, 4).Value = 552.6299192227849

$ws = $wb.Worksheets.Item("8")
$ws.Cells.Item(2, 4).Value = 0.10618158520092773
$ws.Cells.Item(2, 5).Value = 194.89722
$ws.Cells.Item(3, 2).Value = -1723.387446692379
$ws.Cells.Item(3, 4).Value = 709.1755555494638

$ws = $wb.Worksheets.Item("9")
$ws.Cells.Item(2, 4).Value = 0.13408658150585936
$ws.Cells.Item(2, 5).Value = 194.01011
$ws.Cells.Item(3, 2).Value = -1733.6804531563703
$ws.Cells.Item(3, 3).Value = 0.07333183421091723
$ws.Cells.Item(3, 4).Value = 617.3412489941653

$ws = $wb.Worksheets.Item("10")
$ws.Cells.Item(2, 4).Value = 0.07894882409387206
$ws.Cells.Item(2, 5).Value = 188.63784
$ws.Cells.Item(3, 2).Value = -1746.9903591802204
$ws.Cells.Item(3, 3).Value = 0.01941701336526609
$ws.Cells.Item(3, 4).Value = 962.779364893814
